$wb = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item("ANDALUCIA")
$ws = $wb.Worksheets.Add($null, $ws0)
$ws.Name = "VALENCIA"
$ws.Range("B2").Value = "Gandía"
$ws.Range("B3").Value = "Xàtiva"
$ws.Range("B4").Value = "L'Olleria"
